$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "278.09"
Set-TextValue "G2" "15"
Set-TextValue "D3" "22.91"
Set-TextValue "G3" "15"
Set-TextValue "D4" "6.373"
Set-TextValue "G4" "15"
Set-TextValue "D5" "0.06265"
Set-TextValue "G5" "15"
Set-TextValue "G6" "15"
Set-TextValue "D7" "6.628"
Set-TextValue "G7" "15"
Set-TextValue "D8" "1.404"
Set-TextValue "G8" "15"
Set-TextValue "D9" "0.8302"
Set-TextValue "G9" "15"
Set-TextValue "G10" "15"
Set-TextValue "D11" "0.1608"
Set-TextValue "G11" "15"
Set-TextValue "D12" "0.08446"
Set-TextValue "G12" "15"
Set-TextValue "D13" "0.03507"
Set-TextValue "G13" "15"
Set-TextValue "D14" "0.03229"
Set-TextValue "G14" "15"
Set-TextValue "D15" "4.077"
Set-TextValue "G15" "15"
Set-TextValue "D16" "0.09288"
Set-TextValue "G16" "15"
Set-TextValue "D17" "0.001668"
Set-TextValue "G17" "15"
Set-TextValue "D18" "0.04753"
Set-TextValue "G18" "15"
Set-TextValue "D19" "0.006340"
Set-TextValue "G19" "15"
Set-TextValue "D20" "0.005732"
Set-TextValue "G20" "15"
Set-TextValue "G21" "15"
Set-TextValue "G22" "15"
Set-TextValue "D23" "3.729"
Set-TextValue "G23" "15"
Set-TextValue "D24" "2.414"
Set-TextValue "G24" "15"
Set-TextValue "D25" "0.3331"
Set-TextValue "G25" "15"
Set-TextValue "D26" "0.1260"
Set-TextValue "G26" "15"
Set-TextValue "G27" "15"
Set-TextValue "D28" "0.0002700"
Set-TextValue "G28" "15"
Set-TextValue "G29" "15"
Set-TextValue "G30" "15"
Set-TextValue "G31" "15"
Set-TextValue "G32" "15"
Set-TextValue "G33" "15"
Set-TextValue "G34" "15"
Set-TextValue "G35" "15"
Set-TextValue "G36" "15"
Set-TextValue "G37" "15"
Set-TextValue "G38" "15"
Set-TextValue "G39" "15"
Set-TextValue "D40" "0.04744"
Set-TextValue "G40" "15"
Set-TextValue "D41" "0.007121"
Set-TextValue "G41" "15"
Set-TextValue "D42" "0.1171"
Set-TextValue "G42" "15"
Set-TextValue "D43" "0.003648"
Set-TextValue "G43" "15"
Set-TextValue "D44" "0.01230"
Set-TextValue "G44" "15"
Set-TextValue "D45" "0.00006092"
Set-TextValue "G45" "15"
Set-TextValue "D46" "0.0009885"
Set-TextValue "G46" "15"
Set-TextValue "G47" "15"
Set-TextValue "D48" "0.7809"
Set-TextValue "G48" "15"
Set-TextValue "G49" "15"
Set-TextValue "D50" "0.00001398"
Set-TextValue "G50" "15"
Set-TextValue "D51" "0.01238"
Set-TextValue "G51" "15"
